# Fix code style issue: use parentheses instead of backslash for line continuation
# Slide 18, shape "TextBox 3" contains the PPO training code snippet.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# Run 28: "loss = pg_loss - ... entropy_loss \" -> "loss = (pg_loss - ... entropy_loss"
$tr.Runs(28).Text = "                loss = (pg_loss - self.config.ent_coef * entropy_loss"

# Run 29: "       + self.config.vf_coef * v_loss" -> "        + self.config.vf_coef * v_loss)"
$tr.Runs(29).Text = "                        + self.config.vf_coef * v_loss)"

# Editing the run text re-triggers the text box's auto-fit layout pass; restore
# the shape's original (unchanged) height so only the two text runs differ.
$shape.Height = 381.6
